# Applies scheduled-runner profit refresh to Asura_Profits sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3366.1667
$ws.Range("I62").Value = 3339.4
$ws.Range("J62").Value = 3500
$ws.Range("K62").Value = 3339.4
$ws.Range("L62").Value = 3500
$ws.Range("M62").Value = -2715.4
$ws.Range("N62").Value = -4748

$ws.Range("H65").Value = 3366.1667
$ws.Range("I65").Value = 3339.4
$ws.Range("J65").Value = 3500
$ws.Range("K65").Value = 16697
$ws.Range("L65").Value = 17500
$ws.Range("M65").Value = -13577
$ws.Range("N65").Value = -23740

$ws.Range("H112").Value = 6546.5186
$ws.Range("J112").Value = 7256.5
$ws.Range("L112").Value = 21769.5
$ws.Range("N112").Value = -23985.5

$ws.Range("H116").Value = 4654354.5
$ws.Range("I116").Value = 5717434.5
$ws.Range("J116").Value = 3380.25
$ws.Range("K116").Value = 5717434.5
$ws.Range("L116").Value = 3380.25
$ws.Range("M116").Value = -5713992.5
$ws.Range("N116").Value = -10264.25

$ws.Range("H129").Value = 1153.7561
$ws.Range("J129").Value = 1166.7693
$ws.Range("L129").Value = 3500.3079
$ws.Range("N129").Value = -13500.3079

$ws.Range("H138").Value = 9533509
$ws.Range("I138").Value = 22227032
$ws.Range("J138").Value = 13366.833
$ws.Range("K138").Value = 66681096
$ws.Range("L138").Value = 40100.499
$ws.Range("M138").Value = -66675956
$ws.Range("N138").Value = -50380.499

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 935.8889
$ws.Range("I45").Value = 637.1667
$ws.Range("J45").Value = 1533.3334
$ws.Range("K45").Value = 637.1667
$ws.Range("L45").Value = 1533.3334
$ws.Range("M45").Value = -260.1667
$ws.Range("N45").Value = -2287.3334

$ws.Range("H135").Value = 59800
$ws.Range("J135").Value = 59800
$ws.Range("L135").Value = 59800
$ws.Range("N135").Value = -69940

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 76299.21000000001
$ws.Range("I20").Value = 114366.445
$ws.Range("J20").Value = 7778.2
$ws.Range("K20").Value = 114366.445
$ws.Range("L20").Value = 7778.2
$ws.Range("M20").Value = -114119.445
$ws.Range("N20").Value = -8272.200000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5438.273
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 5438.273
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 5438.273
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -6028.273

$ws.Range("H34").Value = 5438.273
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 5438.273
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 5438.273
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -5842.273

$ws.Range("H105").Value = 6410.1055
$ws.Range("I105").Value = 7842.2856
$ws.Range("J105").Value = 2400
$ws.Range("K105").Value = 7842.2856
$ws.Range("L105").Value = 2400
$ws.Range("M105").Value = -6095.2856
$ws.Range("N105").Value = -5894

$ws.Range("H135").Value = 124761.664
$ws.Range("J135").Value = 124761.664
$ws.Range("L135").Value = 124761.664
$ws.Range("N135").Value = -134901.664

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 263500
$ws.Range("J15").Value = 263500
$ws.Range("L15").Value = 263500
$ws.Range("N15").Value = -264076

$ws.Range("H80").Value = 3300.6
$ws.Range("J80").Value = 4002
$ws.Range("L80").Value = 4002
$ws.Range("N80").Value = -5998

$ws.Range("H81").Value = 263500
$ws.Range("J81").Value = 263500
$ws.Range("L81").Value = 263500
$ws.Range("N81").Value = -265496

$ws.Range("H83").Value = 3300.6
$ws.Range("J83").Value = 4002
$ws.Range("L83").Value = 20010
$ws.Range("N83").Value = -29994

$ws.Range("H84").Value = 263500
$ws.Range("J84").Value = 263500
$ws.Range("L84").Value = 790500
$ws.Range("N84").Value = -800484

$ws.Range("H102").Value = 1768.1578
$ws.Range("I102").Value = 1750.2307
$ws.Range("J102").Value = 1807
$ws.Range("K102").Value = 1750.2307
$ws.Range("L102").Value = 1807
$ws.Range("M102").Value = -128.2307000000001
$ws.Range("N102").Value = -5051

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4603.5713
$ws.Range("I7").Value = 4585
$ws.Range("J7").Value = 4650
$ws.Range("K7").Value = 4585
$ws.Range("L7").Value = 4650
$ws.Range("M7").Value = -4473
$ws.Range("N7").Value = -4874

$ws.Range("H120").Value = 46523.5
$ws.Range("J120").Value = 46523.5
$ws.Range("L120").Value = 46523.5
$ws.Range("N120").Value = -56199.5

$ws.Range("H122").Value = 13981348
$ws.Range("I122").Value = 19964290
$ws.Range("J122").Value = 21150.834
$ws.Range("K122").Value = 59892870
$ws.Range("L122").Value = 63452.50199999999
$ws.Range("M122").Value = -59890420
$ws.Range("N122").Value = -68352.50199999999

$ws.Range("H126").Value = 4603.5713
$ws.Range("I126").Value = 4585
$ws.Range("J126").Value = 4650
$ws.Range("K126").Value = 13755
$ws.Range("L126").Value = 13950
$ws.Range("M126").Value = -11285
$ws.Range("N126").Value = -18890

$ws.Range("H135").Value = 175065.8
$ws.Range("J135").Value = 175065.8
$ws.Range("L135").Value = 175065.8
$ws.Range("N135").Value = -185205.8

$ws.Range("H138").Value = 65120.832
$ws.Range("J138").Value = 65120.832
$ws.Range("L138").Value = 65120.832
$ws.Range("N138").Value = -75400.83199999999

$ws.Range("H140").Value = 84595
$ws.Range("J140").Value = 84595
$ws.Range("L140").Value = 84595
$ws.Range("N140").Value = -94955

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 44910
$ws.Range("J16").Value = 44910
$ws.Range("L16").Value = 44910
$ws.Range("N16").Value = -45494

$ws.Range("H96").Value = 897.6667
$ws.Range("I96").Value = 846.5
$ws.Range("K96").Value = 846.5
$ws.Range("M96").Value = 526.5

$ws.Range("H120").Value = 35259.8
$ws.Range("J120").Value = 35259.8
$ws.Range("L120").Value = 35259.8
$ws.Range("N120").Value = -44935.8

$ws.Range("H121").Value = 27884.076
$ws.Range("J121").Value = 27884.076
$ws.Range("L121").Value = 27884.076
$ws.Range("N121").Value = -31378.076

$ws.Range("H126").Value = 8203.727999999999
$ws.Range("I126").Value = 10693.25
$ws.Range("J126").Value = 1565
$ws.Range("K126").Value = 32079.75
$ws.Range("L126").Value = 4695
$ws.Range("M126").Value = -29609.75
$ws.Range("N126").Value = -9635
